# Rename the embedded picture id/name labels (wp:docPr/@name and
# pic:cNvPr/@name) on the three logo images living in the first-page
# footer, the default footer, and the first-page header:
#   footer (first page) id=3 : image2.png -> image1.png
#   footer (default)    id=2 : image2.png -> image1.png
#   header (first page) id=1 : image1.jpg -> image2.jpg
#
# InlineShape has no settable "Name" property (matches real Word), so
# the only way to touch wp:docPr/@name is to delete the picture and
# re-insert an (otherwise identical) drawing XML fragment with the
# corrected name at the same spot.

$d = $word.ActiveDocument
$sec = $d.Sections.Item(1)

function Replace-LogoName($range, [string]$oldRunXml, [string]$newName) {
    $shape = $range.InlineShapes.Item(1)
    $ip = $shape.Range
    $shape.Delete()

    $xml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships">
<w:body>
<w:p>
$oldRunXml
</w:p>
</w:body>
</w:document>
</pkg:xmlData></pkg:part>
</pkg:package>
"@

    $ip.InsertXML($xml)
}

# --- Footer, first page (footer1.xml, docPr id="3"): image2.png -> image1.png ---
$footerFirst = $sec.Footers.Item(2)
$runFooterFirst = '<w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Verdana" w:cs="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana"/><w:b w:val="1"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="none"/><w:shd w:fill="auto" w:val="clear"/><w:vertAlign w:val="baseline"/></w:rPr><w:drawing><wp:inline distB="0" distT="0" distL="0" distR="0"><wp:extent cx="952500" cy="285750"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="3" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId1"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="952500" cy="285750"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>'
Replace-LogoName $footerFirst.Range $runFooterFirst "image1.png"

# --- Footer, default (footer2.xml, docPr id="2"): image2.png -> image1.png ---
$footerDefault = $sec.Footers.Item(1)
$runFooterDefault = '<w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:b w:val="1"/></w:rPr><w:drawing><wp:inline distB="0" distT="0" distL="0" distR="0"><wp:extent cx="952500" cy="285750"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="2" name="image1.png"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="Y:\Together Design\Pearson Edexcel PowerPoint amends\Assets\PearsonLogo.png" id="0" name="image1.png"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId1"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="952500" cy="285750"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>'
Replace-LogoName $footerDefault.Range $runFooterDefault "image1.png"

# --- Header, first page (header1.xml, docPr id="1"): image1.jpg -> image2.jpg ---
$headerFirst = $sec.Headers.Item(2)
$runHeaderFirst = '<w:r w:rsidDel="00000000" w:rsidR="00000000" w:rsidRPr="00000000"><w:rPr><w:rFonts w:ascii="Verdana" w:cs="Verdana" w:eastAsia="Verdana" w:hAnsi="Verdana"/><w:b w:val="0"/><w:i w:val="0"/><w:smallCaps w:val="0"/><w:strike w:val="0"/><w:color w:val="000000"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:u w:val="none"/><w:shd w:fill="auto" w:val="clear"/><w:vertAlign w:val="baseline"/></w:rPr><w:drawing><wp:inline distB="0" distT="0" distL="0" distR="0"><wp:extent cx="914400" cy="277792"/><wp:effectExtent b="0" l="0" r="0" t="0"/><wp:docPr descr="BTec_Logo-Orange" id="1" name="image2.jpg"/><a:graphic><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic><pic:nvPicPr><pic:cNvPr descr="BTec_Logo-Orange" id="0" name="image2.jpg"/><pic:cNvPicPr preferRelativeResize="0"/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId1"/><a:srcRect b="0" l="0" r="0" t="0"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="914400" cy="277792"/></a:xfrm><a:prstGeom prst="rect"/><a:ln/></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r>'
Replace-LogoName $headerFirst.Range $runHeaderFirst "image2.jpg"

Write-Host "Renamed 3 inline picture labels (footer id=3, footer id=2, header id=1)."
